# Add four new prediction columns (WIN, TOP2, TOP4, RELEGATION) between the
# existing "Team" and "ExpPoints" columns, shifting "ExpPoints" from C to G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank columns at C:F (shifts old column C "ExpPoints" -> G),
# preserving the header style that was on C1.
$ws.Range("C1:F1").EntireColumn.Insert()

# Populate the new header cells.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"
